$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1730.2858
$ws.Cells.Item(19, 9).Value = 1179
$ws.Cells.Item(19, 11).Value = 1179
$ws.Cells.Item(19, 13).Value = -1004

$ws.Cells.Item(28, 8).Value = 405404.9
$ws.Cells.Item(28, 9).Value = 667276.0600000001
$ws.Cells.Item(28, 11).Value = 667276.0600000001
$ws.Cells.Item(28, 13).Value = -666791.0600000001

$ws.Cells.Item(33, 8).Value = 382.5
$ws.Cells.Item(33, 9).Value = 267
$ws.Cells.Item(33, 10).Value = 729
$ws.Cells.Item(33, 11).Value = 267
$ws.Cells.Item(33, 12).Value = 729
$ws.Cells.Item(33, 13).Value = -38
$ws.Cells.Item(33, 14).Value = -1187

$ws.Cells.Item(38, 8).Value = 2443.625
$ws.Cells.Item(38, 9).Value = 254.5
$ws.Cells.Item(38, 10).Value = 4632.75
$ws.Cells.Item(38, 11).Value = 763.5
$ws.Cells.Item(38, 12).Value = 13898.25
$ws.Cells.Item(38, 13).Value = -391.5
$ws.Cells.Item(38, 14).Value = -14642.25

$ws.Cells.Item(39, 8).Value = 598.8461
$ws.Cells.Item(39, 10).Value = 3499.5
$ws.Cells.Item(39, 12).Value = 10498.5
$ws.Cells.Item(39, 14).Value = -11090.5

$ws.Cells.Item(40, 8).Value = 2689.75
$ws.Cells.Item(40, 10).Value = 2379
$ws.Cells.Item(40, 12).Value = 2379
$ws.Cells.Item(40, 14).Value = -2729

$ws.Cells.Item(41, 8).Value = 335.77777
$ws.Cells.Item(41, 9).Value = 315.375
$ws.Cells.Item(41, 10).Value = 499
$ws.Cells.Item(41, 11).Value = 315.375
$ws.Cells.Item(41, 12).Value = 499
$ws.Cells.Item(41, 13).Value = 124.625
$ws.Cells.Item(41, 14).Value = -1379

$ws.Cells.Item(42, 8).Value = 152.66667
$ws.Cells.Item(42, 9).Value = 58
$ws.Cells.Item(42, 11).Value = 174
$ws.Cells.Item(42, 13).Value = 56

$ws.Cells.Item(64, 8).Value = 16682056
$ws.Cells.Item(64, 9).Value = 47630188
$ws.Cells.Item(64, 11).Value = 47630188
$ws.Cells.Item(64, 13).Value = -47629940

$ws.Cells.Item(67, 8).Value = 16682056
$ws.Cells.Item(67, 9).Value = 47630188
$ws.Cells.Item(67, 11).Value = 47630188
$ws.Cells.Item(67, 13).Value = -47629330

$ws.Cells.Item(82, 8).Value = 6042
$ws.Cells.Item(82, 10).Value = 16615.666
$ws.Cells.Item(82, 12).Value = 49846.99800000001
$ws.Cells.Item(82, 14).Value = -50658.99800000001

$ws.Cells.Item(85, 8).Value = 6042
$ws.Cells.Item(85, 10).Value = 16615.666
$ws.Cells.Item(85, 12).Value = 49846.99800000001
$ws.Cells.Item(85, 14).Value = -52654.99800000001

$ws.Cells.Item(86, 8).Value = 1430382.1
$ws.Cells.Item(86, 9).Value = 2501475.5
$ws.Cells.Item(86, 10).Value = 2257.6667
$ws.Cells.Item(86, 11).Value = 2501475.5
$ws.Cells.Item(86, 12).Value = 2257.6667
$ws.Cells.Item(86, 13).Value = -2500352.5
$ws.Cells.Item(86, 14).Value = -4503.6667

$ws.Cells.Item(89, 8).Value = 1430382.1
$ws.Cells.Item(89, 9).Value = 2501475.5
$ws.Cells.Item(89, 10).Value = 2257.6667
$ws.Cells.Item(89, 11).Value = 12507377.5
$ws.Cells.Item(89, 12).Value = 11288.3335
$ws.Cells.Item(89, 13).Value = -12501761.5
$ws.Cells.Item(89, 14).Value = -22520.3335

$ws.Cells.Item(100, 8).Value = 7294.381
$ws.Cells.Item(100, 10).Value = 11699.9
$ws.Cells.Item(100, 12).Value = 11699.9
$ws.Cells.Item(100, 14).Value = -12781.9

$ws.Cells.Item(103, 8).Value = 547.3333
$ws.Cells.Item(103, 10).Value = 548.25
$ws.Cells.Item(103, 12).Value = 1644.75
$ws.Cells.Item(103, 14).Value = -2816.75

$ws.Cells.Item(115, 8).Value = 2431.4285
$ws.Cells.Item(115, 10).Value = 3185
$ws.Cells.Item(115, 12).Value = 9555
$ws.Cells.Item(115, 14).Value = -12689

$ws.Cells.Item(131, 8).Value = 1510.6
$ws.Cells.Item(131, 9).Value = 919.5
$ws.Cells.Item(131, 10).Value = 3875
$ws.Cells.Item(131, 11).Value = 2758.5
$ws.Cells.Item(131, 12).Value = 11625
$ws.Cells.Item(131, 13).Value = 2281.5
$ws.Cells.Item(131, 14).Value = -21705

$ws.Cells.Item(137, 8).Value = 7181.2583
$ws.Cells.Item(137, 9).Value = 2420.1428
$ws.Cells.Item(137, 10).Value = 11102.177
$ws.Cells.Item(137, 11).Value = 7260.428400000001
$ws.Cells.Item(137, 12).Value = 33306.531
$ws.Cells.Item(137, 13).Value = -4710.428400000001
$ws.Cells.Item(137, 14).Value = -38406.531

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 2958.2942
$ws.Cells.Item(74, 9).Value = 2160.389
$ws.Cells.Item(74, 11).Value = 2160.389
$ws.Cells.Item(74, 13).Value = -1286.389

$ws.Cells.Item(77, 8).Value = 2958.2942
$ws.Cells.Item(77, 9).Value = 2160.389
$ws.Cells.Item(77, 11).Value = 10801.945
$ws.Cells.Item(77, 13).Value = -6433.945

$ws.Cells.Item(97, 8).Value = 403.10526
$ws.Cells.Item(97, 9).Value = 347.15384
$ws.Cells.Item(97, 10).Value = 524.3333
$ws.Cells.Item(97, 11).Value = 347.15384
$ws.Cells.Item(97, 12).Value = 524.3333
$ws.Cells.Item(97, 13).Value = 148.84616
$ws.Cells.Item(97, 14).Value = -1516.3333

$ws.Cells.Item(102, 8).Value = 2902.1304
$ws.Cells.Item(102, 9).Value = 2865.4211
$ws.Cells.Item(102, 10).Value = 3076.5
$ws.Cells.Item(102, 11).Value = 2865.4211
$ws.Cells.Item(102, 12).Value = 3076.5
$ws.Cells.Item(102, 13).Value = -1243.4211
$ws.Cells.Item(102, 14).Value = -6320.5

$ws.Cells.Item(122, 8).Value = 3345.9148
$ws.Cells.Item(122, 9).Value = 3054.7334
$ws.Cells.Item(122, 11).Value = 9164.200199999999
$ws.Cells.Item(122, 13).Value = -6714.200199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4998.533
$ws.Cells.Item(86, 9).Value = 1373.75
$ws.Cells.Item(86, 11).Value = 1373.75
$ws.Cells.Item(86, 13).Value = -250.75

$ws.Cells.Item(89, 8).Value = 4998.533
$ws.Cells.Item(89, 9).Value = 1373.75
$ws.Cells.Item(89, 11).Value = 6868.75
$ws.Cells.Item(89, 13).Value = -1252.75

$ws.Cells.Item(99, 8).Value = 4518.5186
$ws.Cells.Item(99, 9).Value = 3490.25
$ws.Cells.Item(99, 10).Value = 7456.4287
$ws.Cells.Item(99, 11).Value = 3490.25
$ws.Cells.Item(99, 12).Value = 7456.4287
$ws.Cells.Item(99, 13).Value = -1992.25
$ws.Cells.Item(99, 14).Value = -10452.4287

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3637.5476
$ws.Cells.Item(31, 10).Value = 7835.4375
$ws.Cells.Item(31, 12).Value = 7835.4375
$ws.Cells.Item(31, 14).Value = -8425.4375

$ws.Cells.Item(34, 8).Value = 3637.5476
$ws.Cells.Item(34, 10).Value = 7835.4375
$ws.Cells.Item(34, 12).Value = 7835.4375
$ws.Cells.Item(34, 14).Value = -8239.4375

$ws.Cells.Item(62, 8).Value = 4902.273
$ws.Cells.Item(62, 10).Value = 6002.467
$ws.Cells.Item(62, 12).Value = 6002.467
$ws.Cells.Item(62, 14).Value = -7250.467

$ws.Cells.Item(65, 8).Value = 4902.273
$ws.Cells.Item(65, 10).Value = 6002.467
$ws.Cells.Item(65, 12).Value = 30012.335
$ws.Cells.Item(65, 14).Value = -36252.335

$ws.Cells.Item(105, 8).Value = 22224434
$ws.Cells.Item(105, 9).Value = 26317410
$ws.Cells.Item(105, 11).Value = 26317410
$ws.Cells.Item(105, 13).Value = -26315663

$ws.Cells.Item(134, 8).Value = 16302260
$ws.Cells.Item(134, 9).Value = 18016734
$ws.Cells.Item(134, 11).Value = 54050202
$ws.Cells.Item(134, 13).Value = -54047667

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 4330.631
$ws.Cells.Item(107, 9).Value = 346
$ws.Cells.Item(107, 10).Value = 5424.451
$ws.Cells.Item(107, 11).Value = 1038
$ws.Cells.Item(107, 12).Value = 16273.353
$ws.Cells.Item(107, 13).Value = 882
$ws.Cells.Item(107, 14).Value = -20113.353

$ws.Cells.Item(118, 8).Value = 14491.2
$ws.Cells.Item(118, 9).Value = 24029
$ws.Cells.Item(118, 10).Value = 12106.75
$ws.Cells.Item(118, 11).Value = 72087
$ws.Cells.Item(118, 12).Value = 36320.25
$ws.Cells.Item(118, 13).Value = -70844
$ws.Cells.Item(118, 14).Value = -38806.25

$ws.Cells.Item(120, 8).Value = 22811.545
$ws.Cells.Item(120, 9).Value = 10132.429
$ws.Cells.Item(120, 11).Value = 30397.287
$ws.Cells.Item(120, 13).Value = -25559.287

$ws.Cells.Item(122, 8).Value = 113930.11
$ws.Cells.Item(122, 9).Value = 289.23077
$ws.Cells.Item(122, 10).Value = 160096.72
$ws.Cells.Item(122, 11).Value = 2603.07693
$ws.Cells.Item(122, 12).Value = 1440870.48
$ws.Cells.Item(122, 13).Value = -153.0769300000002
$ws.Cells.Item(122, 14).Value = -1445770.48

$ws.Cells.Item(131, 8).Value = 41275296
$ws.Cells.Item(131, 9).Value = 53334108
$ws.Cells.Item(131, 11).Value = 160002324
$ws.Cells.Item(131, 13).Value = -159997284

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 983.53845
$ws.Cells.Item(97, 10).Value = 1730.5
$ws.Cells.Item(97, 12).Value = 1730.5
$ws.Cells.Item(97, 14).Value = -2722.5

$ws.Cells.Item(107, 8).Value = 333.5
$ws.Cells.Item(107, 9).Value = 270.22223
$ws.Cells.Item(107, 11).Value = 270.22223
$ws.Cells.Item(107, 13).Value = 1649.77777

$ws.Cells.Item(122, 8).Value = 8989.706
$ws.Cells.Item(122, 9).Value = 9509.1
$ws.Cells.Item(122, 11).Value = 28527.3
$ws.Cells.Item(122, 13).Value = -26077.3

$ws.Cells.Item(126, 10).Value = 10289.947
$ws.Cells.Item(126, 12).Value = 30869.841
$ws.Cells.Item(126, 14).Value = -35809.841

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 19232076
$ws.Cells.Item(46, 9).Value = 1195.238
$ws.Cells.Item(46, 11).Value = 1195.238
$ws.Cells.Item(46, 13).Value = -1007.238

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 15928.286
$ws.Cells.Item(62, 9).Value = 14749.667
$ws.Cells.Item(62, 10).Value = 16812.25
$ws.Cells.Item(62, 11).Value = 14749.667
$ws.Cells.Item(62, 12).Value = 16812.25
$ws.Cells.Item(62, 13).Value = -14125.667
$ws.Cells.Item(62, 14).Value = -18060.25

$ws.Cells.Item(65, 8).Value = 15928.286
$ws.Cells.Item(65, 9).Value = 14749.667
$ws.Cells.Item(65, 10).Value = 16812.25
$ws.Cells.Item(65, 11).Value = 73748.33499999999
$ws.Cells.Item(65, 12).Value = 84061.25
$ws.Cells.Item(65, 13).Value = -70628.33499999999
$ws.Cells.Item(65, 14).Value = -90301.25

$ws.Cells.Item(132, 8).Value = 6800.613
$ws.Cells.Item(132, 9).Value = 5891.45
$ws.Cells.Item(132, 11).Value = 17674.35
$ws.Cells.Item(132, 13).Value = -15144.35
